$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 134; this pushes the existing rows
# 134..212 down to 135..213 (and extends the used range accordingly).
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with a new data record (same
# "template" fields as the surrounding rows, new date/price figures).
$ws.Range("A134").Value = 3
$ws.Range("B134").Value = "Femacal de La Calera"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44488
$ws.Range("E134").Value = 5
$ws.Range("F134").Value = 100112009
$ws.Range("G134").Value = "Acelga"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 270
$ws.Range("K134").Value = 2000
$ws.Range("L134").Value = 2200
$ws.Range("M134").Value = 2089
$ws.Range("N134").Value = '$/docena de atados (6 kilos)'
$ws.Range("O134").Value = "Provincia de Quillota"
$ws.Range("P134").Value = 348
$ws.Range("Q134").Value = 6
$ws.Range("R134").Value = "Hortaliza"
